# Edit script implementing the diff between before.docx and after.docx

$d = $word.ActiveDocument

function Find-ParaIndex($text, $styleName) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $pp = $d.Paragraphs.Item($i)
        $t = $pp.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            if ($styleName -eq $null -or $pp.Style.NameLocal -eq $styleName) {
                return $i
            }
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) Merge the two runs "System " + "developers are expected..." into one run
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "System developers are expected to adhere to these Principles, requiring Governance based Decisions to deviate from them.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "System developers are expected to adhere to these Principles, requiring Governance based Decisions to deviate from them.",
    2)

# ---------------------------------------------------------------------------
# 2) Insert the new "PRINC-DEV-01: SOLID" / "PRINC-DEV-01: Avoid Class
#    Nesting" sections right before the "Appendices" heading.
# ---------------------------------------------------------------------------
$appendicesIdx = Find-ParaIndex "Appendices" "Appendices"
$anchorRange = $d.Paragraphs.Item($appendicesIdx).Range

# Insert 15 empty paragraphs before the Appendices heading; we will fill
# each one in afterwards (content + style).
for ($k = 0; $k -lt 15; $k++) {
    $anchorRange.InsertParagraphBefore()
}

# Recompute the index after insertion - the 15 new (still-empty) paragraphs
# now occupy the slots immediately before the (shifted) Appendices heading.
$appendicesIdx = Find-ParaIndex "Appendices" "Appendices"
$base = $appendicesIdx - 15

# -- paragraph base+0 : empty BodyText --------------------------------------
$p0 = $d.Paragraphs.Item($base + 0)
$p0.Style = "BodyText"

# -- paragraph base+1 : Heading2 "PRINC-DEV-01: SOLID" ----------------------
$p1 = $d.Paragraphs.Item($base + 1)
$p1.Style = "Heading2"
$p1.Range.Text = "PRINC-DEV-01: SOLID"

# -- paragraph base+2 : BodyText ---------------------------------------------
$p2 = $d.Paragraphs.Item($base + 2)
$p2.Style = "BodyText"
$p2.Range.Text = "Develop Object Oriented (O.O.) code following S.O.L.I.D. principles."

# -- paragraph base+3 : Heading3 "Details" -----------------------------------
$p3 = $d.Paragraphs.Item($base + 3)
$p3.Style = "Heading3"
$p3.Range.Text = "Details"

# -- paragraph base+4 : BodyText ---------------------------------------------
$p4 = $d.Paragraphs.Item($base + 4)
$p4.Style = "BodyText"
$p4.Range.Text = "SOLID is an acronym for the first 5 principles (of 10) best practice principles for Object Oriented (O.O.) development:"

# -- paragraphs base+5 .. base+9 : the 5 bullet list items -------------------
$li1 = $d.Paragraphs.Item($base + 5)
$li1.Style = "ListParagraph"
$li1.Range.Text = "single responsibility principle: in essence, code that has a single purpose is smaller, more predictable, while also reducing refactoring and maintenance as removes the need to be changed, unless the specific purpose changes.  "

$li2 = $d.Paragraphs.Item($base + 6)
$li2.Style = "ListParagraph"
$li2.Range.Text = "open-closed principle, "

$li3 = $d.Paragraphs.Item($base + 7)
$li3.Style = "ListParagraph"
$li3.Range.Text = "Liskov substitution principle: in essence reference interfaces. "

$li4 = $d.Paragraphs.Item($base + 8)
$li4.Style = "ListParagraph"
$li4.Range.Text = "interface segregation principle, and "

$li5 = $d.Paragraphs.Item($base + 9)
$li5.Style = "ListParagraph"
$li5.Range.Text = "dependency inversion principle."

# Apply the bullet list to all five items at once (keeps a single numId).
$listStart = $d.Paragraphs.Item($base + 5).Range.Start
$listEnd = $d.Paragraphs.Item($base + 9).Range.End
$listRange = $d.Range($listStart, $listEnd)
$listRange.ListFormat.ApplyBulletDefault()

# -- paragraph base+10 : Heading3 "Recommendations" --------------------------
$p10 = $d.Paragraphs.Item($base + 10)
$p10.Style = "Heading3"
$p10.Range.Text = "Recommendations"

# -- paragraph base+11 : Heading2 "PRINC-DEV-01: Avoid Class Nesting" -------
$p11 = $d.Paragraphs.Item($base + 11)
$p11.Style = "Heading2"
$p11.Range.Text = "PRINC-DEV-01: Avoid Class Nesting"

# -- paragraph base+12 : Heading3 "Details" ----------------------------------
$p12 = $d.Paragraphs.Item($base + 12)
$p12.Style = "Heading3"
$p12.Range.Text = "Details"

# -- paragraph base+13 : Normal body paragraph -------------------------------
$p13 = $d.Paragraphs.Item($base + 13)
$p13.Style = "Normal"
$p13.Range.Text = "Develop classes based on interfaces and abstract base classes, while avoiding to nest base classes to deeply (1 is ok, 3 is suspect).  "

# -- paragraph base+14 : completely empty paragraph --------------------------
$p14 = $d.Paragraphs.Item($base + 14)
$p14.Style = "Normal"

